# ---------------------------------------------------------------------------
# Electricity Technology Shareweights.xlsx - "CA 4.0 files test" edit
#
# Summary of changes applied:
#  1. About sheet: insert 2 rows at row 10, add a new bolded note row
#     ("The EPS assumes shareweights will be between 0 and 1 (inclusive).")
#  2. ETS sheet: split "natural gas nonpeaker" into two technologies
#     ("natural gas steam turbine" in place, "natural gas combined cycle"
#     as a new row right after it), and append 7 brand-new technology rows
#     at the bottom of the table (CCS variants, SMR, hydrogen turbines),
#     plus flip the "municipal solid waste" shareweight from 0 to 1.
#  3. Re-point the active sheet/selection/zoom to match the saved view.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$ets   = $wb.Worksheets.Item("ETS")

# ---------------------------------------------------------------------------
# 1) About sheet
# ---------------------------------------------------------------------------

# Insert two new rows before the old row 10, shifting the remainder of the
# notes paragraph down by two rows (old r10 -> r12, ... old r21 -> r23).
$about.Range("A10:A11").EntireRow.Insert()

$about.Range("A10").Value = "The EPS assumes shareweights will be between 0 and 1 (inclusive)."
$about.Range("A10").Font.Bold = $true

# ---------------------------------------------------------------------------
# 2) ETS sheet - technology table
# ---------------------------------------------------------------------------

# Insert one new row right after row 3 ("natural gas nonpeaker"); this
# shifts every technology from the old row 4 onward down by one row.
$ets.Range("A4:AF4").EntireRow.Insert()

# Row 3 was "natural gas nonpeaker" -> rename in place to the steam-turbine
# technology (shareweight values for this row are unchanged, still all 1).
$ets.Range("A3").Value = "natural gas steam turbine"

# New row 4: "natural gas combined cycle", shareweight = 1 for every year.
$ets.Range("A4").Value = "natural gas combined cycle"
$ets.Range("B4:AF4").Value = 1

# Existing technology rows 5-17 (after the insert) keep their prior
# shareweight values - only the new tech-split above required the insert.
# Row 18 ("municipal solid waste") flips from 0 to 1.
$ets.Range("B18:AF18").Value = 1

# Append seven brand-new technology rows (19-25) with shareweight = 1 for
# every year.
$newTechRows = @(
    @(19, "hard coal w CCS"),
    @(20, "natural gas combined cycle w CCS"),
    @(21, "biomass w CCS"),
    @(22, "lignite w CCS"),
    @(23, "small modular reactor"),
    @(24, "hydrogen combustion turbine"),
    @(25, "hydrogen combined cycle")
)

foreach ($item in $newTechRows) {
    $r = $item[0]
    $name = $item[1]
    $ets.Range("A$r").Value = $name
    $ets.Range("B${r}:AF${r}").Value = 1
}

# The last two new rows (hydrogen technologies) get a distinct font: solid
# black text color, vertically centered. Apply directly to A24, then copy
# the resulting format onto A25 so only a single new font/style pair is
# added to the style table (matches xlPasteFormats = -4122).
$ets.Range("A24").VerticalAlignment = -4108
$ets.Range("A24").Font.Color = 0
$ets.Range("A24").Copy()
$ets.Range("A25").PasteSpecial(-4122)
$ets.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Views / selection / active sheet
# ---------------------------------------------------------------------------

$about.Select()
$about.Range("A25").Select()
$excel.ActiveWindow.Zoom = 130

$ets.Range("B25:AF25").Select()
$ets.Application.ActiveWindow.ScrollRow = 1
$ets.Application.ActiveWindow.ScrollColumn = 5

$about.Select()
